$wb = $excel.ActiveWorkbook

# --- Sheet "21" (top 10 species, area 21): remove "Pagophilus groenlandicus"
# (a seal, mis-included in the species aggregate due to an ISSCAAP grouping
# bug) from the top-10 ranking. Every subsequent row moves up one slot and a
# new 10th-place species ("Crassostrea virginica") appears at the bottom.
$ws21 = $wb.Worksheets.Item("21")
$ws21.Range("A7").Value = "Chionoecetes opilio"
$ws21.Range("B7").Value = 80249.2
$ws21.Range("A8").Value = "Clupea harengus"
$ws21.Range("B8").Value = 73585.92
$ws21.Range("A9").Value = "Reinhardtius hippoglossoides"
$ws21.Range("B9").Value = 66899.18000000001
$ws21.Range("A10").Value = "Spisula solidissima"
$ws21.Range("B10").Value = 60674.75
$ws21.Range("A11").Value = "Crassostrea virginica"
$ws21.Range("B11").Value = 47148.78

# --- Sheet "27": "Laminaria hyperborea" (a seaweed, same bug) drops out of
# the 10th spot and is replaced by "Sebastes mentella".
$ws27 = $wb.Worksheets.Item("27")
$ws27.Range("A11").Value = "Sebastes mentella"
$ws27.Range("B11").Value = 94462.19

# --- Sheet "87" (area 87): "Lessonia nigrescens" (a kelp, same bug) drops
# out of the top 10 and every following row moves up one slot; a new
# 10th-place species ("Merluccius gayi") appears at the bottom.
$ws87 = $wb.Worksheets.Item("87")
$ws87.Range("A8").Value = "Thunnus albacares"
$ws87.Range("B8").Value = 117107.95
$ws87.Range("A9").Value = "Actinopterygii"
$ws87.Range("B9").Value = 100798.73
$ws87.Range("A10").Value = "Sarda chiliensis"
$ws87.Range("B10").Value = 94502.89999999999
$ws87.Range("A11").Value = "Merluccius gayi"
$ws87.Range("B11").Value = 87525.53
